# Actualización automática 2025-08-19 17:15:09
#
# Updates sales figures for CASTRO ALCIVAR EDA MARIA across the three
# sheets of the workbook: "VENTAS POR GRUPO" (by client/category),
# "VENTA MENSUAL" (by month) and "CUMPLIMIENTO MENSUAL" (budget vs.
# actual summary), reflecting a new PIEDRA SINTERIZADA sale recorded for
# RIVERA CANTOS MARTHA JACQUELIN (295.63) in agosto, and doubling the
# 240X80 PORCELANATO figure recorded for RUIZ PINEDA LUIS ALFREDO.

$wb = $excel.ActiveWorkbook

# --- Hoja "VENTAS POR GRUPO" ---------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# RIVERA CANTOS MARTHA JACQUELIN ahora registra una venta de PIEDRA SINTERIZADA
$ws1.Range("L43").Value = 295.63

# RUIZ PINEDA LUIS ALFREDO: 240X80 PORCELANATO se duplica
$ws1.Range("D44").Value = 915.84

# Conteo de asesores con venta en PIEDRA SINTERIZADA sube de 3 a 4 (de 55)
$ws1.Range("L57").Value = "4 de 55"

# --- Hoja "VENTA MENSUAL" -------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# agosto: nueva venta para RIVERA CANTOS MARTHA JACQUELIN
$ws2.Range("F43").Value = 295.63

# agosto: venta actualizada para RUIZ PINEDA LUIS ALFREDO
$ws2.Range("F44").Value = 1042.1

# Total de agosto
$ws2.Range("F57").Value = 20394.25

# --- Hoja "CUMPLIMIENTO MENSUAL" -----------------------------------------
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# 240X80 PORCELANATO
$ws3.Range("D3").Value = 15059.81
$ws3.Range("E3").Value = -175.9300000000003
$ws3.Range("F3").Value = 1.011820170546927

# PIEDRA SINTERIZADA
$ws3.Range("D15").Value = 4486.16
$ws3.Range("E15").Value = 16203.84
$ws3.Range("F15").Value = 0.2168274528757854

# PORCELANATO
$ws3.Range("D16").Value = 12154.53
$ws3.Range("E16").Value = 46566.7
$ws3.Range("F16").Value = 0.2069869789852835

# TOTAL
$ws3.Range("D19").Value = 33345.41
$ws3.Range("E19").Value = 76523.34
$ws3.Range("F19").Value = 0.3035022242448376

# La columna E (VENTA) se ensancha ligeramente para el nuevo valor
$ws3.Columns.Item(5).ColumnWidth = 24 - 0.8333333333333334
